$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("O").Delete() | Out-Null

$ws.Range("O4").Select() | Out-Null
